# Work Time Table edit
# - Set cell B5 to "6 timer"
# - Scroll the sheet view so that topLeftCell is A4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the value of B5 (adds a new shared string "6 timer")
$ws.Range("B5").Value = "6 timer"

# Select B5 (keep existing selection) and scroll the view so A4 is the top-left visible cell
$ws.Activate()
$ws.Range("B5").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
